$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" value would be auto-parsed as a number by Excel
# need to be forced to text format first so the exact display string
# (including trailing zeros / leading zeros) is preserved, matching the
# inline-string content from the source diff.

$ws.Range("D2").Value = '27.434.92'
$ws.Range("E2").Value = '  +2.44%  '

$ws.Range("D3").Value = '1.799.45'
$ws.Range("E3").Value = '  +3.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '338.58'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3808'
$ws.Range("E7").Value = '  +1.51%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3459'
$ws.Range("E8").Value = '  +2.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.73'
$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.205'
$ws.Range("E10").Value = '  +1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07531'
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.13'
$ws.Range("E13").Value = '  +8.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.494'
$ws.Range("E14").Value = '  +1.44%  '

$ws.Range("D15").Value = '1.797.21'
$ws.Range("E15").Value = '  +3.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.100'
$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001108'
$ws.Range("E17").Value = '  +2.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06658'
$ws.Range("E18").Value = '  -1.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9996'
$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.539'
$ws.Range("E21").Value = '  +4.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.42'
$ws.Range("E22").Value = '  +4.04%  '

$ws.Range("D23").Value = '27.428.15'
$ws.Range("E23").Value = '  +2.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.60'
$ws.Range("E24").Value = '  -1.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("E25").Value = '  -1.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.572'
$ws.Range("E26").Value = '  +6.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.499'
$ws.Range("E27").Value = '  +1.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.52'
$ws.Range("E28").Value = '  +9.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.34'
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").Value = '2.001.64'
$ws.Range("E30").Value = '  +3.36%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '133.92'
$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.053'
$ws.Range("E32").Value = '  -1.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.141'
$ws.Range("E33").Value = '  +1.84%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08708'
$ws.Range("E34").Value = '  +0.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.33'
$ws.Range("E35").Value = '  +3.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.646'
$ws.Range("E36").Value = '  -2.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.467'
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6922'
$ws.Range("E38").Value = '  +10.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.921'
$ws.Range("E39").Value = '  +5.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06403'
$ws.Range("E40").Value = '  +2.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2209'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02345'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.275'
$ws.Range("E43").Value = '  +4.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.46'
$ws.Range("E44").Value = '  +1.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6469'
$ws.Range("E45").Value = '  +6.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.869'
$ws.Range("E47").Value = '  -1.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.145'
$ws.Range("E48").Value = '  +3.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '130.62'
$ws.Range("E49").Value = '  +1.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07204'
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.08'
$ws.Range("E51").Value = '  +2.78%  '
